# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Valor Mora" (F column) amounts for the periods 1910 (row 16) and
# 1807 (row 22) were swapped in the source database update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = 44000
$ws.Range("F22").Value = 42000
